$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "29.869.36"
$c.Style = "Normal"
$ws.Range("E2").Value = "  -0.03%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.639.72"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +1.05%  "
$ws.Range("E4").Value = "  +0.67%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "215.45"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.93%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.519"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +0.16%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +0.63%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "28.72"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -2.55%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.260"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +0.62%  "
$ws.Range("E10").Value = "  +0.22%  "
$ws.Range("E11").Value = "  -1.15%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "1.875.03"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +1.02%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "1.639.73"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +1.15%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "0.591"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +4.14%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "9.44"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +7.01%  "
$ws.Range("E16").Value = "  -1.55%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "29.882.62"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -0.08%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "64.61"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +0.48%  "
$ws.Range("E19").Value = "  -0.82%  "
$ws.Range("E20").Value = "  -0.56%  "
$ws.Range("E21").Value = "  +0.57%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "9.90"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +3.24%  "
$ws.Range("E23").Value = "  +0.91%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "2.19"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +2.48%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "157.50"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +0.53%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "15.54"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -0.39%  "
$ws.Range("E27").Value = "  -0.90%  "
$ws.Range("E29").Value = "  +0.57%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "0.0493"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +1.20%  "
$ws.Range("E31").Value = "  -0.49%  "
$ws.Range("E32").Value = "  +1.70%  "
$ws.Range("E33").Value = "  -0.68%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "1.425.03"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +0.06%  "
$ws.Range("E35").Value = "  +3.24%  "
$ws.Range("E36").Value = "  -0.74%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "2.73"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -4.62%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.0174"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +2.20%  "
$ws.Range("E39").Value = "  +0.08%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "76.59"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +10.75%  "
$ws.Range("E41").Value = "  +0.87%  "
$ws.Range("E42").Value = "  +0.94%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.0499"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -1.32%  "
$ws.Range("E44").Value = "  -0.40%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +0.64%  "
$ws.Range("E46").Value = "  -1.00%  "
$ws.Range("E47").Value = "  -0.68%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "1.781.37"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +0.92%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "49.13"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -8.94%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "93.48"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +6.08%  "
$ws.Range("E51").Value = "  -1.65%  "
